$d = $word.ActiveDocument

# --- Header info block ---

# "AVANCE DE PROYECTO N° ___" -> "AVANCE DE PROYECTO: EV-PARCIAL"
$d.Content.Find.Execute("AVANCE DE PROYECTO N° ___", $true, $false, $false, $false, $false, `
    $true, 1, $false, "AVANCE DE PROYECTO: EV-PARCIAL", 2) | Out-Null

# "CURSO: " -> "CURSO: Diseño y Desarrollo Web"
$d.Content.Find.Execute("CURSO: ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "CURSO: Diseño y Desarrollo Web", 2) | Out-Null

# "NRC:" -> "NRC: 1806"
$d.Content.Find.Execute("NRC:", $true, $false, $false, $false, $false, `
    $true, 1, $false, "NRC: 1806", 2) | Out-Null

# "FECHA DE PRESENTACIÓN: " -> "FECHA DE PRESENTACIÓN: 23/10/2024"
$d.Content.Find.Execute("FECHA DE PRESENTACIÓN: ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "FECHA DE PRESENTACIÓN: 23/10/2024", 2) | Out-Null

# "SEMANA LECTIVA: " -> "SEMANA LECTIVA: Sem-08"
$d.Content.Find.Execute("SEMANA LECTIVA: ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "SEMANA LECTIVA: Sem-08", 2) | Out-Null

# --- Table: "Cuadro de aportes individuales" ---

$t = $d.Tables.Item(1)

# Row 3 (member 1): name after "1."
$d.Content.Find.Execute("1.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "1. Pinedo Gutierrez Christopher David", 2) | Out-Null

# Row 3 (member 1): aporte description (empty cell, col 2)
$cell = $t.Cell(3, 2)
$cell.Range.Text = "Creación de guión y desarrollo de la parte final (Animaciones y despedida.)"
$t.Cell(3, 2).Range.Font.Name = "Arial"

# EXCELENTE column (col 3) marked with "x" for rows 3-7 (members 1-5)
for ($row = 3; $row -le 7; $row++) {
    $cell = $t.Cell($row, 3)
    $cell.Range.Text = "x"
    $t.Cell($row, 3).Range.Font.Name = "Arial"
}

Write-Host "done"
